# Logboek.xlsx update — add a new logged activity row, and move the
# YouTube-link footer cell down to make (visual) room, per the author's
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the hyperlink footer cell from L31 to L30 -------------------
# Cut/paste the cell (keeps the L31 placeholder's style behind) then
# re-point the hyperlink itself at the new location, reusing the
# worksheet's existing "Hyperlink" cell style so no duplicate style is
# created.
$ws.Range("L31").Cut($ws.Range("L30"))
$excel.CutCopyMode = 0

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("L30"), "https://www.youtube.com/watch?v=z0MimkXIvE8") | Out-Null
$ws.Range("L30").Style = "Hyperlink"

# --- 2. Fill in the new logboek entry on row 19 ---------------------------
# Copy the banded row formatting (fill + wrap + date format) down from the
# row above's matching style band so B19/C19/D19 pick up the same styles
# Excel would apply when typing into an already-formatted table row.
$ws.Range("B17:D17").Copy()
$ws.Range("B19:D19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D19").Value = "180 minuten"
$ws.Range("B19").Value = "Extra postman videos bekeken, Environment en collection runner uitgewerkt"
$ws.Range("C19").Value = 44537
$ws.Rows.Item(19).RowHeight = 45

# --- 3. Last selected cell -------------------------------------------------
$ws.Range("M35").Select()
